# Update "想去人数" (interested-count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 686
$ws1.Range("F3").Value = 527
$ws1.Range("F7").Value = 48
$ws1.Range("F8").Value = 3345
$ws1.Range("F9").Value = 4274
$ws1.Range("F10").Value = 122

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 64

# --- Sheet "全部类型" (all types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 686
$ws4.Range("F3").Value = 527
$ws4.Range("F7").Value = 48
$ws4.Range("F8").Value = 3345
$ws4.Range("F9").Value = 4274
$ws4.Range("F10").Value = 122
$ws4.Range("F11").Value = 64
